$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F "想去人数" (number interested) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F11").Value = 1348
$wsExpo.Range("F12").Value = 3007
$wsExpo.Range("F13").Value = 428
$wsExpo.Range("F14").Value = 1628
$wsExpo.Range("F17").Value = 239
$wsExpo.Range("F18").Value = 1391
$wsExpo.Range("F24").Value = 3490
$wsExpo.Range("F25").Value = 689
$wsExpo.Range("F27").Value = 1546

# Sheet "本地生活" (Local life) - column F update
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 7

# Sheet "全部类型" (All types) - column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 7
$wsAll.Range("F21").Value = 1348
$wsAll.Range("F22").Value = 3007
$wsAll.Range("F23").Value = 428
$wsAll.Range("F24").Value = 1628
$wsAll.Range("F27").Value = 239
$wsAll.Range("F28").Value = 1391
$wsAll.Range("F36").Value = 3490
$wsAll.Range("F37").Value = 689
$wsAll.Range("F39").Value = 1546
